$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (e.g. "68.193.29",
# "0.999", "10.90"). The source workbook stores them as text (inline
# strings), not numbers, and some use "." as a thousands separator or
# have significant trailing zeros. Force these cells to Text format before
# writing so Excel does not auto-convert/round them to numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.193.29"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "3.659.99"
$ws.Range("E3").Value = "  -4.33%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "591.81"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "180.37"
$ws.Range("E6").Value = "  +8.63%  "
$ws.Range("D7").Value = "3.654.30"
$ws.Range("E7").Value = "  -4.43%  "
$ws.Range("D8").Value = "0.628"
$ws.Range("E8").Value = "  -5.03%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "0.713"
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -7.00%  "
$ws.Range("D12").Value = "56.13"
$ws.Range("E12").Value = "  +7.02%  "
$ws.Range("D13").Value = "0.0000292"
$ws.Range("E13").Value = "  -8.00%  "
$ws.Range("E14").Value = "  -5.22%  "
$ws.Range("D15").Value = "4.230.70"
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("D16").Value = "3.650.66"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").Value = "19.28"
$ws.Range("E17").Value = "  -6.45%  "
$ws.Range("E18").Value = "  -2.01%  "
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  -6.75%  "
$ws.Range("D20").Value = "1.11"
$ws.Range("E20").Value = "  -6.61%  "
$ws.Range("D21").Value = "67.803.78"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "408.99"
$ws.Range("E22").Value = "  -5.45%  "
$ws.Range("D23").Value = "4.56"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("D24").Value = "88.46"
$ws.Range("E24").Value = "  -4.76%  "
$ws.Range("D25").Value = "4.11"
$ws.Range("E25").Value = "  +3.74%  "
$ws.Range("D26").Value = "3.01"
$ws.Range("E26").Value = "  -7.43%  "
$ws.Range("D27").Value = "12.73"
$ws.Range("E27").Value = "  -7.39%  "
$ws.Range("D28").Value = "10.90"
$ws.Range("E28").Value = "  -4.74%  "
$ws.Range("D29").Value = "6.04"
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "9.47"
$ws.Range("E30").Value = "  -9.10%  "
$ws.Range("D31").Value = "32.57"
$ws.Range("E31").Value = "  -6.13%  "
$ws.Range("D32").Value = "7.21"
$ws.Range("E32").Value = "  -11.02%  "
$ws.Range("D33").Value = "12.35"
$ws.Range("E33").Value = "  -7.63%  "
$ws.Range("D34").Value = "0.117"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("D35").Value = "64.68"
$ws.Range("E35").Value = "  -5.48%  "
$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "43.22"
$ws.Range("E36").Value = "  -8.78%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "600.39"
$ws.Range("E37").Value = "  -6.11%  "
$ws.Range("D38").Value = "0.0₃0891"
$ws.Range("E38").Value = "  -9.55%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.398"
$ws.Range("E40").Value = "  -7.05%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  -6.40%  "
$ws.Range("D43").Value = "3.01"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("E44").Value = "  -7.12%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0436"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  -12.17%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "0.134"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "9.01"
$ws.Range("E48").Value = "  -8.25%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.70"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.723.31"
$ws.Range("E50").Value = "  -4.49%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "3.13"
$ws.Range("E51").Value = "  -2.98%  "
